# Add season-record columns (Wins / Losses / Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the formatting of the last existing header cell (AC1,
# which is bold/centered/bordered) onto the three new header cells, then
# set their text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows: every player on the roster shares the team's 1998 season
# record - 92 wins, 70 losses, 0 ties.
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 92
    $ws.Cells.Item($r, 31).Value = 70
    $ws.Cells.Item($r, 32).Value = 0
}
